# Update row 2 values per target diff:
#  A2: 75c44810a32a3d6447df      -> 7f065c23251bf386d439
#  B2: +74267426006              -> +74267426011   (must stay text, keep leading "+")
#  C2: Automation User 10        -> Automation User 11
#  E2: 2025-12-30                -> 2025-12-31      (must stay text, not become a date serial)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the cells whose new values would otherwise be
# auto-converted by Excel (phone number -> number, date string -> date).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "7f065c23251bf386d439"
$ws.Range("B2").Value = "+74267426011"
$ws.Range("C2").Value = "Automation User 11"
$ws.Range("E2").Value = "2025-12-31"
